# Fill in the previously-placeholder value for 2025-11-04 (row 102, B column)
# and append a new row 103 for 2025-11-05 on each of the 3 sheets.

$wb = $excel.ActiveWorkbook

# New row values per sheet (in sheet order: 한화솔루션, 아난티, 대아티아이)
$newRow102B = @(488151, 59787, 13642)
$newRow103A = @(45966, 45966, 45966)
$newRow103B = @(442298, 61329, 13625)

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)

    # Update the B102 value (was a placeholder 0)
    $ws.Range("B102").Value = $newRow102B[$i - 1]

    # Add the new row 103
    $ws.Range("A103").Value = $newRow103A[$i - 1]
    $ws.Range("A103").NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Range("B103").Value = $newRow103B[$i - 1]
}
